$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.034.23"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "1.863.13"
$ws.Range("E3").Value = "  -2.26%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "307.11"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5128"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.78%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3744"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07120"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8876"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.57"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07566"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.864.75"
$ws.Range("E13").Value = "  -2.24%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.315"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "88.69"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.65%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008477"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.67%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.12"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.12%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "27.029.83"
$ws.Range("E20").Value = "  -2.93%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.050"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").Value = "2.106.69"
$ws.Range("E22").Value = "  -0.78%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.54"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.72%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.470"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.83%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "149.76"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.94%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.846"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.46%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.96"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.14%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.099"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.27%  "
$ws.Range("E29").Value = "  -1.63%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.684"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.07%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.658"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09035"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.74%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05126"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.05%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.081"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.33%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.153"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -6.44%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7329"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.86%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02048"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.508"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.81%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.058"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.073"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.66%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5312"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.583"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "115.69"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.17%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.284"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.08%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1469"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.78%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4612"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("E48").Value = "  -5.20%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.567"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.99%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "64.28"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "36.61"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
